# Auto-generated PowerShell Excel COM-interop script
# Applies the cryptos.xlsx price/volume/ranking update described in the commit diff
# (GitHub Actions scheduled refresh, Wed Nov 22 15:52:59 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.514.09'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '2.031.19'
$ws.Range("E3").Value = '  +2.19%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.10'
$ws.Range("E5").Value = '  -12.59%  '
$ws.Range("E6").Value = '  -1.79%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.14'
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.14'
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0745'
$ws.Range("E11").Value = '  -1.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.100'
$ws.Range("E12").Value = '  -0.15%  '
$ws.Range("D13").Value = '2.329.63'
$ws.Range("E13").Value = '  +3.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.27'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.03'
$ws.Range("E15").Value = '  -6.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.760'
$ws.Range("E16").Value = '  -1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.11'
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").Value = '2.033.78'
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("D19").Value = '36.749.98'
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.41'
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.60'
$ws.Range("E21").Value = '  +10.84%  '
$ws.Range("D22").Value = '0.0₃0794'
$ws.Range("E22").Value = '  -3.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '220.61'
$ws.Range("E23").Value = '  -5.69%  '
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.38'
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("E26").Value = '  -8.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.66'
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.65'
$ws.Range("E28").Value = '  -2.01%  '
$ws.Range("E29").Value = '  +3.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.88'
$ws.Range("E30").Value = '  -2.07%  '
$ws.Range("E31").Value = '  +1.37%  '
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.36'
$ws.Range("E33").Value = '  -3.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0602'
$ws.Range("E34").Value = '  -3.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.50'
$ws.Range("E35").Value = '  +4.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.26'
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("E38").Value = '  -2.55%  '
$ws.Range("E39").Value = '  +9.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.27'
$ws.Range("E40").Value = '  -5.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.92'
$ws.Range("E41").Value = '  -2.15%  '
$ws.Range("D42").Value = '1.471.61'
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0930'
$ws.Range("E43").Value = '  +1.88%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '92.60'
$ws.Range("E44").Value = '  +5.08%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.15'
$ws.Range("E45").Value = '  +36.63%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0203'
$ws.Range("E46").Value = '  -1.30%  '
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.12'
$ws.Range("E47").Value = '  -4.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.64'
$ws.Range("E48").Value = '  +1.33%  '
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.89'
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.90'
$ws.Range("E51").Value = '  +1.47%  '
